$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add 3 new columns (K, L, M) ---
$ws.Range('J1').Copy()
$ws.Range('K1:M1').PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range('K1').Value = 'diff_lightgbm'
$ws.Range('L1').Value = 'diff_catboost'
$ws.Range('M1').Value = 'diff_randomforest'

# --- Data rows 2-16: rewritten in new order, plus 3 new diff columns ---
$data = @(
    @('Cytarabine', 'OC[C@H]1O[C@H]([C@H]([C@@H]1O)O)n1ccc(nc1=O)N', -5.299697499, 'weak', 305, -3.680272961593789, -3.273539318661022, -3.935420296470931, -3.629744192241914, 1.669953306758086, 1.364277202529069, 2.026158180338977, 1.619424537406211),
    @('N4-acetyl-4-amino-1-[(2R,3S,4R,5R)-3,4-dihydroxy-5 (hydroxymethyl)oxolan-2-yl] pyrimidin-2-one', 'OC[C@H]1O[C@H]([C@H]([C@H]1O)O)n1ccc(nc1=O)NC(=O)C', -3.332697499, 'weak', 305, -3.554533138279848, -3.183357587695099, -3.944206091934551, -3.560698939303166, 0.2280014403031658, 0.6115085929345505, 0.1493399113049012, 0.2218356392798482),
    @('N4-hexanoyl-4-amino-1-[(2R,3S,4R,5R)-3,4-dihydroxy-5 (hydroxymethyl)oxolan-2-yl] pyrimidin-2-one', 'CCCCCC(=O)Nc1ccn(c(=O)n1)[C@@H]1O[C@@H]([C@@H]([C@@H]1O)O)CO', -2.180697499, 'weak', 305, -3.020732738095238, -2.967456893680463, -3.231940255559159, -3.07337662911162, 0.8926791301116199, 1.051242756559159, 0.7867593946804634, 0.8400352390952377),
    @('Flurbiprofen', 'OC(=O)C(c1ccc(c(c1)F)c1ccccc1)C', -1.774697499, 'ionized', 310, -2.937572831686134, -2.663225373077253, -2.801034395641547, -2.800610866801645, 1.025913367801645, 1.026336896641547, 0.8885278740772533, 1.162875332686135),
    @('Flurbiprofen glucoside', 'OC1O[C@H](COC(=O)C(c2ccc(c(c2)F)c2ccccc2)C)[C@H]([C@@H]([C@H]1O)O)O', -2.920697499, 'unionizable', 310, -2.808845093590896, -2.572122776855409, -2.591462968254274, -2.657476946233526, 0.2632205527664739, 0.3292345307457261, 0.3485747221445914, 0.1118524054091043),
    @('Flurbiprofen mannoside', 'OC1O[C@H](COC(=O)C(c2ccc(c(c2)F)c2ccccc2)C)[C@H]([C@@H]([C@@H]1O)O)O', -3.378697499, 'unionizable', 310, -2.808845093590896, -2.572122776855409, -2.591462968254274, -2.657476946233526, 0.7212205527664737, 0.7872345307457258, 0.8065747221445911, 0.5698524054091041),
    @('Ibuprofen', 'CC(Cc1ccc(cc1)C(C(=O)O)C)C', -1.946697499, 'ionized', 310, -2.936437117400419, -2.625305516140017, -2.716059308418287, -2.759267313986241, 0.8125698149862408, 0.7693618094182872, 0.6786080171400168, 0.9897396184004188),
    @('Ibuprofen glucoside', 'CC(Cc1ccc(cc1)C(C(=O)OC[C@H]1OC(O)[C@@H]([C@H]([C@@H]1O)O)O)C)C', -3.388697499, 'unionizable', 310, -2.835328426924228, -2.503273457045836, -2.591462968254274, -2.643354950741446, 0.7453425482585541, 0.797234530745726, 0.8854240419541641, 0.5533690720757716),
    @('Ibuprofen mannoside', 'CC(Cc1ccc(cc1)C(C(=O)OC[C@H]1OC(O)[C@H]([C@H]([C@@H]1O)O)O)C)C', -2.991697499, 'unionizable', 310, -2.835328426924228, -2.503273457045836, -2.591462968254274, -2.643354950741446, 0.3483425482585538, 0.4002345307457258, 0.4884240419541639, 0.1563690720757713),
    @('Ketoprofen', 'OC(=O)C(c1cccc(c1)C(=O)c1ccccc1)C', -2.594697499, 'ionized', 310, -2.937572831686134, -2.661001120813594, -2.787188265324236, -2.795254072607988, 0.2005565736079884, 0.1924907663242359, 0.0663036218135944, 0.3428753326861345),
    @('Ketoprofen glucoside', 'OC1O[C@H](COC(=O)C(c2cccc(c2)C(=O)c2ccccc2)C)[C@H]([C@@H]([C@H]1O)O)O', -4.174697499, 'unionizable', 310, -2.819765093590896, -2.562996992615822, -2.591462968254274, -2.658075018153664, 1.516622480846336, 1.583234530745726, 1.611700506384177, 1.354932405409103),
    @('Ketoprofen mannoside', 'OC1O[C@H](COC(=O)C(c2cccc(c2)C(=O)c2ccccc2)C)[C@H]([C@@H]([C@@H]1O)O)O', -4.037697499, 'unionizable', 310, -2.819765093590896, -2.562996992615822, -2.591462968254274, -2.658075018153664, 1.379622480846336, 1.446234530745726, 1.474700506384178, 1.217932405409104),
    @('Naproxen', 'COc1ccc2c(c1)ccc(c2)C(C(=O)O)C', -3.661697499, 'ionized', 310, -2.883814991964879, -2.602661556308178, -2.801034395641547, -2.762503647971535, 0.8991938510284654, 0.8606631033584531, 1.059035942691823, 0.777882507035121),
    @('Naproxen glucoside', 'COc1ccc2c(c1)ccc(c2)[C@@H](C(=O)OC[C@H]1OC(O)[C@@H]([C@H]([C@@H]1O)O)O)C', -4.329697499, 'unionizable', 310, -2.797832593590896, -2.563592053913609, -2.518280394512903, -2.626568347339136, 1.703129151660864, 1.811417104487097, 1.76610544508639, 1.531864905409104),
    @('Naproxen mannoside', 'COc1ccc2c(c1)ccc(c2)[C@@H](C(=O)OC[C@H]1OC(O)[C@H]([C@H]([C@@H]1O)O)O)C', -4.143697499, 'unionizable', 310, -2.797832593590896, -2.563592053913609, -2.518280394512903, -2.626568347339136, 1.517129151660864, 1.625417104487097, 1.580105445086391, 1.345864905409104)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# --- Sheet-level structure updates ---
$ws.Columns.Item(2).ColumnWidth = 8.140625
$ws.Columns.Item(11).ColumnWidth = 13.140625
$ws.Columns.Item(12).ColumnWidth = 12.7109375
$ws.Columns.Item(13).ColumnWidth = 17.5703125

# Reset AutoFilter to cover the new range A1:M16 (also drops old sort state)
$ws.AutoFilterMode = $false
$ws.Range('A1:M16').AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range
$n = $wb.Names.Item(1)
$n.RefersTo = '=Sheet1!\$A\$1:\$M\$16'

# Restore selection to match the authored state
$ws.Range('L18').Select()

